$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are forced to text (leading apostrophe) so numeric-looking
# strings such as "1.000" or "0.5000" are not coerced into plain numbers,
# matching the workbook's original inline-string text representation.

$ws.Range("D2").Value = "'30.921.24"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'246.72"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5000"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").Value = "'0.3000"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.06863"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "'1.907.84"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").Value = "'17.39"
$ws.Range("E11").Value = "  +3.57%  "
$ws.Range("D12").Value = "'0.07351"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "'92.22"
$ws.Range("E13").Value = "  +7.30%  "
$ws.Range("D14").Value = "'5.127"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("D15").Value = "'0.6851"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "'30.904.13"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "'0.000008103"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "'13.48"
$ws.Range("E18").Value = "  +5.90%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'2.153.51"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'0.9991"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'4.890"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").Value = "'182.51"
$ws.Range("E23").Value = "  +35.17%  "
$ws.Range("D24").Value = "'6.112"
$ws.Range("E24").Value = "  +9.13%  "
$ws.Range("D25").Value = "'9.404"
$ws.Range("E25").Value = "  +3.06%  "
$ws.Range("D26").Value = "'154.39"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").Value = "'18.83"
$ws.Range("E27").Value = "  +12.26%  "
$ws.Range("D28").Value = "'1.958"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").Value = "'1.396"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'4.395"
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").Value = "'0.08995"
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("D32").Value = "'4.084"
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("D33").Value = "'0.05329"
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("D34").Value = "'0.7540"
$ws.Range("E34").Value = "  +7.11%  "
$ws.Range("D35").Value = "'1.148"
$ws.Range("E35").Value = "  +3.83%  "
$ws.Range("D36").Value = "'2.697"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").Value = "'0.01926"
$ws.Range("E37").Value = "  +16.95%  "
$ws.Range("D38").Value = "'2.732"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'2.204"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").Value = "'0.9418"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'0.4412"
$ws.Range("E41").Value = "  +5.29%  "
$ws.Range("D42").Value = "'106.48"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").Value = "'5.882"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'7.809"
$ws.Range("E45").Value = "  +4.87%  "
$ws.Range("D46").Value = "'0.1370"
$ws.Range("E46").Value = "  +9.04%  "
$ws.Range("D47").Value = "'0.05852"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  +6.12%  "
$ws.Range("D49").Value = "'8.632"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "'1.399"
$ws.Range("E51").Value = "  +4.61%  "
